# Adds rows 17-21 to the "Artfynd" sheet: five new species-observation
# records (Goodyera repens / Knärot sightings at Brunnsjöberget, Dlr,
# reported 2023-08-31 by Philipp Weiss), extending the sheet dimension
# from A1:AY16 to A1:AY21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # The diff stores Startdatum/Starttid/Slutdatum/Sluttid as plain text
    # (e.g. "2023-08-31", "00:00"), but a direct .Value assignment of such
    # a string gets auto-parsed into a real Excel date/time serial. Route
    # it through a text-literal formula, then paste-special "values only"
    # over itself so the final cell holds a plain string with no formula
    # and no extra number-format style left behind.
    $c = $ws.Cells.Item($row, $col)
    $escaped = $text -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $ws.Calculate()
    $c.Copy()
    $c.PasteSpecial(-4163) # xlPasteValues
    $ws.Application.CutCopyMode = $false
}

# ----- Row 17 -----
$ws.Cells.Item(17, 1).Value = 111821926
$ws.Cells.Item(17, 2).Value = 96348
$ws.Cells.Item(17, 3).Value = "Ovaliderad"
$ws.Cells.Item(17, 4).Value = "VU"
$ws.Cells.Item(17, 5).Value = 220787
$ws.Cells.Item(17, 6).Value = "Knärot"
$ws.Cells.Item(17, 7).Value = "Goodyera repens"
$ws.Cells.Item(17, 8).Value = "(L.) R. Br."
$ws.Cells.Item(17, 16).Value = "Brunnsjöberget, Dlr"
$ws.Cells.Item(17, 17).Value = 550846.2444635418
$ws.Cells.Item(17, 18).Value = 6681625.195240833
$ws.Cells.Item(17, 19).Value = 15
$ws.Cells.Item(17, 20).Value = "Dalarna"
$ws.Cells.Item(17, 21).Value = "Hedemora"
$ws.Cells.Item(17, 22).Value = "Dalarna"
$ws.Cells.Item(17, 23).Value = "Hedemora"
Set-TextValue 17 25 "2023-08-31"
Set-TextValue 17 26 "00:00"
Set-TextValue 17 27 "2023-08-31"
Set-TextValue 17 28 "00:00"
$ws.Cells.Item(17, 30).Value = $false
$ws.Cells.Item(17, 31).Value = $false
$ws.Cells.Item(17, 33).Value = $false
$ws.Cells.Item(17, 49).Value = "Philipp Weiss"
$ws.Cells.Item(17, 50).Value = "Philipp Weiss"

# ----- Row 18 -----
$ws.Cells.Item(18, 1).Value = 111821927
$ws.Cells.Item(18, 2).Value = 96348
$ws.Cells.Item(18, 3).Value = "Ovaliderad"
$ws.Cells.Item(18, 4).Value = "VU"
$ws.Cells.Item(18, 5).Value = 220787
$ws.Cells.Item(18, 6).Value = "Knärot"
$ws.Cells.Item(18, 7).Value = "Goodyera repens"
$ws.Cells.Item(18, 8).Value = "(L.) R. Br."
$ws.Cells.Item(18, 16).Value = "Brunnsjöberget, Dlr"
$ws.Cells.Item(18, 17).Value = 550819.8901872271
$ws.Cells.Item(18, 18).Value = 6681733.007140613
$ws.Cells.Item(18, 19).Value = 15
$ws.Cells.Item(18, 20).Value = "Dalarna"
$ws.Cells.Item(18, 21).Value = "Hedemora"
$ws.Cells.Item(18, 22).Value = "Dalarna"
$ws.Cells.Item(18, 23).Value = "Hedemora"
Set-TextValue 18 25 "2023-08-31"
Set-TextValue 18 26 "00:00"
Set-TextValue 18 27 "2023-08-31"
Set-TextValue 18 28 "00:00"
$ws.Cells.Item(18, 30).Value = $false
$ws.Cells.Item(18, 31).Value = $false
$ws.Cells.Item(18, 33).Value = $false
$ws.Cells.Item(18, 49).Value = "Philipp Weiss"
$ws.Cells.Item(18, 50).Value = "Philipp Weiss"

# ----- Row 19 -----
$ws.Cells.Item(19, 1).Value = 111821924
$ws.Cells.Item(19, 2).Value = 96348
$ws.Cells.Item(19, 3).Value = "Ovaliderad"
$ws.Cells.Item(19, 4).Value = "VU"
$ws.Cells.Item(19, 5).Value = 220787
$ws.Cells.Item(19, 6).Value = "Knärot"
$ws.Cells.Item(19, 7).Value = "Goodyera repens"
$ws.Cells.Item(19, 8).Value = "(L.) R. Br."
$ws.Cells.Item(19, 11).Value = "blomning"
$ws.Cells.Item(19, 16).Value = "Brunnsjöberget, Dlr"
$ws.Cells.Item(19, 17).Value = 550675.3931295178
$ws.Cells.Item(19, 18).Value = 6681937.422269406
$ws.Cells.Item(19, 19).Value = 15
$ws.Cells.Item(19, 20).Value = "Dalarna"
$ws.Cells.Item(19, 21).Value = "Hedemora"
$ws.Cells.Item(19, 22).Value = "Dalarna"
$ws.Cells.Item(19, 23).Value = "Hedemora"
Set-TextValue 19 25 "2023-08-31"
Set-TextValue 19 26 "00:00"
Set-TextValue 19 27 "2023-08-31"
Set-TextValue 19 28 "00:00"
$ws.Cells.Item(19, 30).Value = $false
$ws.Cells.Item(19, 31).Value = $false
$ws.Cells.Item(19, 33).Value = $false
$ws.Cells.Item(19, 49).Value = "Philipp Weiss"
$ws.Cells.Item(19, 50).Value = "Philipp Weiss"

# ----- Row 20 -----
$ws.Cells.Item(20, 1).Value = 111821928
$ws.Cells.Item(20, 2).Value = 96348
$ws.Cells.Item(20, 3).Value = "Ovaliderad"
$ws.Cells.Item(20, 4).Value = "VU"
$ws.Cells.Item(20, 5).Value = 220787
$ws.Cells.Item(20, 6).Value = "Knärot"
$ws.Cells.Item(20, 7).Value = "Goodyera repens"
$ws.Cells.Item(20, 8).Value = "(L.) R. Br."
$ws.Cells.Item(20, 16).Value = "Brunnsjöberget, Dlr"
$ws.Cells.Item(20, 17).Value = 550825.9503372401
$ws.Cells.Item(20, 18).Value = 6681726.144349095
$ws.Cells.Item(20, 19).Value = 15
$ws.Cells.Item(20, 20).Value = "Dalarna"
$ws.Cells.Item(20, 21).Value = "Hedemora"
$ws.Cells.Item(20, 22).Value = "Dalarna"
$ws.Cells.Item(20, 23).Value = "Hedemora"
Set-TextValue 20 25 "2023-08-31"
Set-TextValue 20 26 "00:00"
Set-TextValue 20 27 "2023-08-31"
Set-TextValue 20 28 "00:00"
$ws.Cells.Item(20, 30).Value = $false
$ws.Cells.Item(20, 31).Value = $false
$ws.Cells.Item(20, 33).Value = $false
$ws.Cells.Item(20, 49).Value = "Philipp Weiss"
$ws.Cells.Item(20, 50).Value = "Philipp Weiss"

# ----- Row 21 -----
$ws.Cells.Item(21, 1).Value = 111821923
$ws.Cells.Item(21, 2).Value = 96348
$ws.Cells.Item(21, 3).Value = "Ovaliderad"
$ws.Cells.Item(21, 4).Value = "VU"
$ws.Cells.Item(21, 5).Value = 220787
$ws.Cells.Item(21, 6).Value = "Knärot"
$ws.Cells.Item(21, 7).Value = "Goodyera repens"
$ws.Cells.Item(21, 8).Value = "(L.) R. Br."
$ws.Cells.Item(21, 11).Value = "blomning"
$ws.Cells.Item(21, 16).Value = "Brunnsjöberget, Dlr"
$ws.Cells.Item(21, 17).Value = 550701.1291094749
$ws.Cells.Item(21, 18).Value = 6681909.496304798
$ws.Cells.Item(21, 19).Value = 15
$ws.Cells.Item(21, 20).Value = "Dalarna"
$ws.Cells.Item(21, 21).Value = "Hedemora"
$ws.Cells.Item(21, 22).Value = "Dalarna"
$ws.Cells.Item(21, 23).Value = "Hedemora"
Set-TextValue 21 25 "2023-08-31"
Set-TextValue 21 26 "00:00"
Set-TextValue 21 27 "2023-08-31"
Set-TextValue 21 28 "00:00"
$ws.Cells.Item(21, 30).Value = $false
$ws.Cells.Item(21, 31).Value = $false
$ws.Cells.Item(21, 33).Value = $false
$ws.Cells.Item(21, 49).Value = "Philipp Weiss"
$ws.Cells.Item(21, 50).Value = "Philipp Weiss"
